# Refresh the cryptos price/volume table (GitHub Actions data pull).
# D = Price, E = Volume(1h); both columns hold text, not numbers, in the
# source sheet. Cells whose new Price value happens to look like a plain
# number (e.g. "0.999", "682.90") are pre-formatted as Text ("@") so the
# Range.Value assignment doesn't get auto-coerced into a numeric cell by
# Excel's type inference - matching the original inlineStr text cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.355.50"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "3.690.24"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "682.90"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.66"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.69%  "
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("E10").Value = "  -3.55%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("E12").Value = "  -2.66%  "
$ws.Range("D13").Value = "4.311.30"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.50"
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("D15").Value = "3.690.33"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "69.336.64"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.09"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.48"
$ws.Range("E19").Value = "  -1.44%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "469.89"
$ws.Range("E20").Value = "  -2.18%  "
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.655"
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "80.04"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").Value = "3.836.13"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  -4.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.95"
$ws.Range("E27").Value = "  -3.68%  "
$ws.Range("E28").Value = "  -3.67%  "
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.75"
$ws.Range("E30").Value = "  -4.11%  "
$ws.Range("E31").Value = "  -2.83%  "
$ws.Range("E32").Value = "  -3.32%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.96"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "3.677.72"
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("E36").Value = "  -6.57%  "
$ws.Range("E37").Value = "  -2.21%  "
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  -4.13%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0909"
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "169.97"
$ws.Range("E43").Value = "  +4.01%  "
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.66"
$ws.Range("E45").Value = "  -0.88%  "
$ws.Range("B46").Value = "SuiNetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.12"
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.72"
$ws.Range("E47").Value = "  -2.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.23"
$ws.Range("E48").Value = "  -6.22%  "
$ws.Range("E49").Value = "  -1.85%  "
$ws.Range("E50").Value = "  -3.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.80"
$ws.Range("E51").Value = "  -3.31%  "
